$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.066.45'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.562.40'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.38'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').Value = '3.023.13'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '63.002.31'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').Value = '2.608.19'
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.48'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +2.73%  '
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  -3.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').Value = '2.693.22'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.61%  '
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('E32').Value = '  +6.86%  '
$ws.Range('D33').Value = '0.0₃0823'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '462.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '175.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.14%  '
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '150.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0546'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.613'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0974'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('E51').Value = '  +0.57%  '
